$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "Bentleigh"
$ws.Cells.Item(2,2).Value = "Il Centro Deli  5/284/292 Centre Rd, Bentleigh VIC 3204"
$ws.Cells.Item(2,3).Value = "22/12/20 12:00pm-12:30pm"
$ws.Cells.Item(2,4).Value = "Case shopped in store"

$ws.Cells.Item(3,1).Value = "Black Rock"
$ws.Cells.Item(3,2).Value = "Woolworths Metro  40 Bluff Road, Black Rock VIC 3193"
$ws.Cells.Item(3,3).Value = "30/12/20 5:30pm-5:55pm"
$ws.Cells.Item(3,4).Value = "Case shopped"

$ws.Cells.Item(4,1).Value = "Box Hill South"
$ws.Cells.Item(4,2).Value = "Bunnings  259 Middleborough Road, Box Hill South VIC 3128"
$ws.Cells.Item(4,3).Value = "30/12/20 12:00pm-12:40pm"
$ws.Cells.Item(4,4).Value = "Case shopped"

$ws.Cells.Item(5,1).Value = "Brighton"
$ws.Cells.Item(5,2).Value = "Brighton Beach  Brighton, VIC 3186"
$ws.Cells.Item(5,3).Value = "29/12/20 12:00pm-3:00pm"
$ws.Cells.Item(5,4).Value = "Case attended beach"

$ws.Cells.Item(6,1).Value = "Camberwell"
$ws.Cells.Item(6,2).Value = "Fu Lin Asian Grocery Supermarket  1397 Toorak Road, Camberwell VIC 3124"
$ws.Cells.Item(6,3).Value = "30/12/20 2:30pm-2:45pm"
$ws.Cells.Item(6,4).Value = "Case shopped"

$ws.Cells.Item(7,1).Value = "Cape Schank"
$ws.Cells.Item(7,2).Value = "National Golf Club  The Cups Drive, Cape Schanck VIC 3939"
$ws.Cells.Item(7,3).Value = "30/12/20 11.40am-1.40pm"
$ws.Cells.Item(7,4).Value = "Case attended course"

$ws.Cells.Item(8,1).Value = "Cheltenham"
$ws.Cells.Item(8,2).Value = "Aldi Cheltenham  280/282 Bay Road, Cheltenham VIC 3192"
$ws.Cells.Item(8,3).Value = "29/12/20 2:00pm-2:30pm"
$ws.Cells.Item(8,4).Value = "Case shopped in store"

$ws.Cells.Item(9,1).Value = "Cheltenham"
$ws.Cells.Item(9,2).Value = "Angus and Cootes Jeweller  Southland Shopping Centre, Shop 2096/1239, Nepean Hwy, Cheltenham VIC 3192"
$ws.Cells.Item(9,3).Value = "28/12/2020 2:30pm-2:50pm"
$ws.Cells.Item(9,4).Value = "Case shopped in store"

$ws.Cells.Item(10,1).Value = "Cheltenham"
$ws.Cells.Item(10,2).Value = "Honey Birdette Southland  Shop 2209/1239, Southland Shopping Centre, Cheltenham VIC 3192"
$ws.Cells.Item(10,3).Value = "22/12/2020 3:50pm-4:05pm"
$ws.Cells.Item(10,4).Value = "Case shopped in store"

$ws.Cells.Item(11,1).Value = "Cheltenham"
$ws.Cells.Item(11,2).Value = "Mecca Southland  Shop 2011/2013, Southland Shopping Centre, Cheltenham VIC 3192"
$ws.Cells.Item(11,3).Value = "22/12/2020 3:30pm-3:50pm"
$ws.Cells.Item(11,4).Value = "Case shopped in store"

$ws.Cells.Item(12,1).Value = "Forest Hill"
$ws.Cells.Item(12,2).Value = "Forest Hill Chase Shopping Centre 270 Canterbury Rd, Forest Hill VIC 3131"
$ws.Cells.Item(12,3).Value = "28/12/20 12:00pm-2:00pm"
$ws.Cells.Item(12,4).Value = "1210hrs Food court 30min; 1250hrs TKMaxx 15min; 1310hrs Target 20min; 1340hrs Woolworths 15min"

$ws.Cells.Item(13,1).Value = "Fountain Gate Shopping Centre"
$ws.Cells.Item(13,2).Value = "Kmart, Big W, Target, Millers, King of Gifts, Lo Costa  25-55 Overland Drive, Narre Warren VIC 3805"
$ws.Cells.Item(13,3).Value = "26/12/20 9:00am-11:00am"
$ws.Cells.Item(13,4).Value = ""

$ws.Cells.Item(14,1).Value = "Glen Waverley"
$ws.Cells.Item(14,2).Value = "Mocha Jos  87 Kingsway, Glen Waverley VIC 3150"
$ws.Cells.Item(14,3).Value = "28/12/20 1:30pm-1:45pm"
$ws.Cells.Item(14,4).Value = ""

$ws.Cells.Item(15,1).Value = "Hallam"
$ws.Cells.Item(15,2).Value = "Coles Hallam  2 Princes Domain Drive, Hallam VIC 3803"
$ws.Cells.Item(15,3).Value = "30/12/20 6:15am-6:30am"
$ws.Cells.Item(15,4).Value = "Case shopped in store"

$ws.Cells.Item(16,1).Value = "Lakes Entrance"
$ws.Cells.Item(16,2).Value = "Blue Riviera Hire Boats  Marine Parade, Lakes Entrance VIC 3909"
$ws.Cells.Item(16,3).Value = "29/12/20 11:15am-12:15pm"
$ws.Cells.Item(16,4).Value = "Case hired a boat"

$ws.Cells.Item(17,1).Value = "Lakes Entrance"
$ws.Cells.Item(17,2).Value = "Central Hotel Lakes Entrance  321 Esplanade, Lakes Entrance VIC 3909"
$ws.Cells.Item(17,3).Value = "30/12/20 5:00pm-6:30pm"
$ws.Cells.Item(17,4).Value = "Case attended outside premises"

$ws.Cells.Item(18,1).Value = "Lakes Entrance"
$ws.Cells.Item(18,2).Value = "Darcey Annas Beach Cafe Kiosk Gift Shop Gallery  426 Main Beach Walk Surf Life Saving, Lakes Entrance VIC 3909"
$ws.Cells.Item(18,3).Value = "30/12/20 11:15am-11:20am"
$ws.Cells.Item(18,4).Value = "Case picked up takeaway"

$ws.Cells.Item(19,1).Value = "Lakes Entrance"
$ws.Cells.Item(19,2).Value = "Woolworths Lakes Entrance 371 Esplanade, Lakes Entrance VIC 3909"
$ws.Cells.Item(19,3).Value = "30/12/20 6:00pm-6:15pm"
$ws.Cells.Item(19,4).Value = "Case shopped in store"

$ws.Cells.Item(20,1).Value = "Mentone"
$ws.Cells.Item(20,2).Value = "Bunnings Mentone  23-27 Nepean Hwy, Mentone VIC 3194"
$ws.Cells.Item(20,3).Value = "29/12/20 07:30am-08:00am"
$ws.Cells.Item(20,4).Value = "Case shopped in store"

$ws.Cells.Item(21,1).Value = "Mentone"
$ws.Cells.Item(21,2).Value = "Bunnings Mentone  23-27 Nepean Hwy, Mentone VIC 3194"
$ws.Cells.Item(21,3).Value = "31/12/20 08:00am-08:30am"
$ws.Cells.Item(21,4).Value = "Case shopped in store"

$ws.Cells.Item(22,1).Value = "Mentone"
$ws.Cells.Item(22,2).Value = "Mentone/Parkdale Beach"
$ws.Cells.Item(22,3).Value = "27/12/20 10:00am-4:30pm"
$ws.Cells.Item(22,4).Value = ""

$ws.Cells.Item(23,1).Value = "Mentone"
$ws.Cells.Item(23,2).Value = "Woolworths Mentone  105-111 Balcombe Road, Mentone VIC 3194"
$ws.Cells.Item(23,3).Value = "23/12/20 2:45pm-3:05pm"
$ws.Cells.Item(23,4).Value = "Case shopped in store"

$ws.Cells.Item(24,1).Value = "Moorabbin"
$ws.Cells.Item(24,2).Value = "COSTCO Moorabbin  8 Chifley Drive, Moorabbin Airport VIC 3194"
$ws.Cells.Item(24,3).Value = "30/12/20 10:45am-12:15pm"
$ws.Cells.Item(24,4).Value = "Case shopped in store"

$ws.Cells.Item(25,1).Value = "Moorabbin"
$ws.Cells.Item(25,2).Value = "COSTCO Moorabbin  8 Chifley Drive, Moorabbin Airport VIC 3194"
$ws.Cells.Item(25,3).Value = "30/12/20 4:00m- 5:50pm"
$ws.Cells.Item(25,4).Value = "Case shopped in store"

$ws.Cells.Item(26,1).Value = "Mordialloc"
$ws.Cells.Item(26,2).Value = "Woodlands Golf Club  109 White Street Mordialloc VIC 3195"
$ws.Cells.Item(26,3).Value = "23/12/20 8:00am-2:00pm"
$ws.Cells.Item(26,4).Value = "Case attended course"

$ws.Cells.Item(27,1).Value = "Mordialloc"
$ws.Cells.Item(27,2).Value = "Woodlands Golf Club  109 White Street Mordialloc VIC 3195"
$ws.Cells.Item(27,3).Value = "28/12/20 12:00pm-6:00pm"
$ws.Cells.Item(27,4).Value = "Case attended course"

$ws.Cells.Item(28,1).Value = "Mount Waverley"
$ws.Cells.Item(28,2).Value = "Ritchies IGA  283 Stephensons Road, Mount Waverley VIC 3149"
$ws.Cells.Item(28,3).Value = "30/12/20 2:00pm-2:30pm"
$ws.Cells.Item(28,4).Value = "Case shopped for half an hour"

$ws.Cells.Item(29,1).Value = "Oakleigh"
$ws.Cells.Item(29,2).Value = "Bunnings Oakleigh  1041 Centre Road, Oakleigh South"
$ws.Cells.Item(29,3).Value = "30/12/20 11:00am-11:30am"
$ws.Cells.Item(29,4).Value = "Case shopped for 30 minutes"

$ws.Cells.Item(30,1).Value = "Oakleigh"
$ws.Cells.Item(30,2).Value = "Katialo restaurant  8 Eaton Mall, Oakleigh VIC 3166"
$ws.Cells.Item(30,3).Value = "28/12/20 7:00pm-7:10pm"
$ws.Cells.Item(30,4).Value = ""

$ws.Cells.Item(31,1).Value = "Wonthaggi"
$ws.Cells.Item(31,2).Value = "Wonthaggi Plaza Shopping centre  2 Biggs Drive, Wonthaggi VIC 3995"
$ws.Cells.Item(31,3).Value = "28/12/20 1:30pm-2.30pm"
$ws.Cells.Item(31,4).Value = "Kmart- shopped for 15 mins"

# Remove the now-obsolete last row (old row 32), shrinking the used range to A1:D31
$ws.Rows.Item(32).Delete()
